$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.588.02"
$ws.Range("E2").Value = "  -2.26%  "
$ws.Range("D3").Value = "1.756.73"
$ws.Range("E3").Value = "  -2.98%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4468"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.64%  "
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07504"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.108"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.060"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.232"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "1.756.87"
$ws.Range("E16").Value = "  -3.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.96"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001065"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06424"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.859"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.03%  "
$ws.Range("D23").Value = "27.628.06"
$ws.Range("E24").Value = "  -2.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.087"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("D28").Value = "1.959.26"
$ws.Range("E28").Value = "  -3.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.137"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.096"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09095"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.571"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.639"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.20"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02305"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2108"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6413"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06016"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.944"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.191"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.000"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.398"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.864"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.80%  "
$ws.Range("E45").Value = "  -3.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5923"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.709"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.977"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "122.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.167"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06872"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.82%  "
